$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "My skill Level" paragraph: insert the new sentences about the helper,
#    and relocate the "_GoBack" bookmark to sit between "...interrupts." and
#    "  Plus, he's a college student..." (matching where the author's cursor
#    was when the edit was last saved).
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("we need done.  I appreciate any help I can get.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    Write-Host "ERROR: could not find skill-level sentence"
}
$start1 = $r1.Start
$rsquo = [string][char]0x2019
$prefixNew = "we need done.  I was getting help from a guy half a world away.  Between the different times, it was difficult to have a conversation.  He" + $rsquo + "s the one that suggested interrupts."
$suffixNew = "  Plus, he's a college student and had very little time to spare.  I appreciate any help I can get."
$r1.Text = $prefixNew + $suffixNew

$bmPos = $start1 + $prefixNew.Length
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 2. Hardware bullet list: "3 PIR motion sensors (1 on my side..." -> "...(One on my side..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("3 PIR motion sensors (1 on my side", $true, $false, $false, $false, $false, $true, 1, $false, "3 PIR motion sensors (One on my side", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Hardware bullet list: `2 "panic" buttons (1 on each bedside table)` -> "(One on each bedside table)"
# ---------------------------------------------------------------------------
$quoteOpen = [string][char]0x201c
$quoteClose = [string][char]0x201d
$panicOld = "2 " + $quoteOpen + "panic" + $quoteClose + " buttons (1 on each bedside table)"
$panicNew = "2 " + $quoteOpen + "panic" + $quoteClose + " buttons (One on each bedside table)"
$d.Content.Find.Execute($panicOld, $true, $false, $false, $false, $false, $true, 1, $false, $panicNew, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Concept bullet list: LDR sentence gains " via PIR's" and " of light ".
# ---------------------------------------------------------------------------
$ldrOld = "LDR prevents lights from turning on during daylight hours.  (This will have to be set to a certain level due to blackout curtains)"
$ldrNew = "LDR prevents lights from turning on via PIR's during daylight hours.  (This will have to be set to a certain level of light due to blackout curtains)"
$d.Content.Find.Execute($ldrOld, $true, $false, $false, $false, $false, $true, 1, $false, $ldrNew, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Concept bullet list: "Fade set by the remote button's code..." -> "...the DIY button's code..."
# ---------------------------------------------------------------------------
$fadeOld = "Fade set by the remote button" + $rsquo + "s code sent via IR emitter)"
$fadeNew = "Fade set by the DIY button" + $rsquo + "s code sent via IR emitter)"
$d.Content.Find.Execute($fadeOld, $true, $false, $false, $false, $false, $true, 1, $false, $fadeNew, 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. Concept bullet list: "G will turn off once power is restored" -> "G will have to turn off once power is restored"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("G will turn off once power is restored", $true, $false, $false, $false, $false, $true, 1, $false, "G will have to turn off once power is restored", 2) | Out-Null

Write-Host "Edits applied"
